$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Update "last updated" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 25 de Marzo de 2020 a las 11:16"

# Simple numeric updates (no reordering): Alemania (row 8), Iran (row 9), Austria (row 15), Brasil (row 22)
Set-Row 8  @("Alemania", 33952, 961, 3299, 30482, 23, 12, 171)
Set-Row 9  @("Iran", 27017, 2206, 8913, 16027, 0, 143, 2077)
Set-Row 15 @("Austria", 5485, 202, 9, 5446, 26, 2, 30)
Set-Row 22 @("Brasil", 2271, 24, 2, 2222, 18, 1, 47)

# Hong Kong's case count overtakes Mexico/Estonia/Egipto/Barein/Argentina, so it moves up
# from row 58 to row 53 (right after Peru), shifting those five countries down by one row.
Set-Row 53 @("Hong Kong", 410, 23, 102, 304, 4, 0, 4)
Set-Row 54 @("Mexico", 405, 38, 4, 396, 1, 1, 5)
Set-Row 55 @("Estonia", 404, 35, 8, 396, 5, 0, 0)
Set-Row 56 @("Egipto", 402, 0, 80, 302, 0, 0, 20)
Set-Row 57 @("Barein", 392, 0, 177, 212, 2, 0, 3)
Set-Row 58 @("Argentina", 387, 0, 52, 329, 0, 0, 6)
